$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 5; this pushes the existing
# rows 5-16 down to rows 7-18 (matching the target diff, which shows the
# old row-5..row-16 data reappearing two rows further down, plus two
# brand-new rows of data at 5 and 6).
$ws.Rows("5:6").Insert()

# --- New row 5 ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 45281
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100103
$ws.Range("H5").Value = "Frutos de hueso (carozo)"
$ws.Range("I5").Value = 100103003
$ws.Range("J5").Value = "Damasco"
$ws.Range("K5").Value = "Castle Brite"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 150
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 28000
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("R5").Value = "Región de Coquimbo"
$ws.Range("S5").Value = 1556
$ws.Range("T5").Value = 18

# --- New row 6 ---
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 45281
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100103
$ws.Range("H6").Value = "Frutos de hueso (carozo)"
$ws.Range("I6").Value = 100103003
$ws.Range("J6").Value = "Damasco"
$ws.Range("K6").Value = "Castle Brite"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 150
$ws.Range("N6").Value = 30000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 30000
$ws.Range("Q6").Value = "$/caja 18 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1667
$ws.Range("T6").Value = 18

# Ensure the Date column keeps its date number format on the two new rows
# (Insert should already have copied format from the row above, but make
# it explicit/robust).
$ws.Range("D5:D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
